# Update TPM-derived NATMI ligand-receptor metrics (columns G:T, rows 2-10)
# to the newly recomputed values (commit: "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.240107999999999
$ws.Range("H2").Value = 18.720324
$ws.Range("I2").Value = 0.01732230523539376
$ws.Range("J2").Value = 0.01732230523539376
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1837056666666667
$ws.Range("N2").Value = 0.5511170000000001
$ws.Range("O2").Value = 0.7269991860920679
$ws.Range("P2").Value = 0.7269991860920678
$ws.Range("Q2").Value = 1.146343200212
$ws.Range("R2").Value = 10.317088801908
$ws.Range("S2").Value = 0.01259330180736963
$ws.Range("T2").Value = 0.01259330180736963
$ws.Range("G3").Value = 6.240107999999999
$ws.Range("H3").Value = 18.720324
$ws.Range("I3").Value = 0.01732230523539376
$ws.Range("J3").Value = 0.01732230523539376
$ws.Range("O3").Value = 0.2534828531892131
$ws.Range("P3").Value = 0.2534828531892131
$ws.Range("Q3").Value = 0.3996955576879999
$ws.Range("R3").Value = 3.597260019191999
$ws.Range("S3").Value = 0.004390907354882054
$ws.Range("T3").Value = 0.004390907354882055
$ws.Range("G4").Value = 6.240107999999999
$ws.Range("H4").Value = 18.720324
$ws.Range("I4").Value = 0.01732230523539376
$ws.Range("J4").Value = 0.01732230523539376
$ws.Range("O4").Value = 0.01951796071871896
$ws.Range("P4").Value = 0.01951796071871896
$ws.Range("Q4").Value = 0.03077621265599999
$ws.Range("R4").Value = 0.276985913904
$ws.Range("S4").Value = 0.0003380960731420752
$ws.Range("T4").Value = 0.0003380960731420752
$ws.Range("I5").Value = 0.9592798330716089
$ws.Range("J5").Value = 0.9592798330716091
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1837056666666667
$ws.Range("N5").Value = 0.5511170000000001
$ws.Range("O5").Value = 0.7269991860920679
$ws.Range("P5").Value = 0.7269991860920678
$ws.Range("Q5").Value = 63.48253877291434
$ws.Range("R5").Value = 571.3428489562291
$ws.Range("S5").Value = 0.6973956578775945
$ws.Range("T5").Value = 0.6973956578775945
$ws.Range("I6").Value = 0.9592798330716089
$ws.Range("J6").Value = 0.9592798330716091
$ws.Range("O6").Value = 0.2534828531892131
$ws.Range("P6").Value = 0.2534828531892131
$ws.Range("S6").Value = 0.2431609890938635
$ws.Range("T6").Value = 0.2431609890938635
$ws.Range("I7").Value = 0.9592798330716089
$ws.Range("J7").Value = 0.9592798330716091
$ws.Range("O7").Value = 0.01951796071871896
$ws.Range("P7").Value = 0.01951796071871896
$ws.Range("S7").Value = 0.01872318610015094
$ws.Range("T7").Value = 0.01872318610015095
$ws.Range("G8").Value = 8.428738666666666
$ws.Range("I8").Value = 0.02339786169299727
$ws.Range("J8").Value = 0.02339786169299728
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.1837056666666667
$ws.Range("N8").Value = 0.5511170000000001
$ws.Range("O8").Value = 0.7269991860920679
$ws.Range("P8").Value = 0.7269991860920678
$ws.Range("Q8").Value = 1.548407055919111
$ws.Range("R8").Value = 13.935663503272
$ws.Range("S8").Value = 0.01701022640710379
$ws.Range("T8").Value = 0.01701022640710379
$ws.Range("G9").Value = 8.428738666666666
$ws.Range("I9").Value = 0.02339786169299727
$ws.Range("J9").Value = 0.02339786169299728
$ws.Range("O9").Value = 0.2534828531892131
$ws.Range("P9").Value = 0.2534828531892131
$ws.Range("Q9").Value = 0.5398831882364443
$ws.Range("S9").Value = 0.00593095674046754
$ws.Range("T9").Value = 0.005930956740467542
$ws.Range("G10").Value = 8.428738666666666
$ws.Range("I10").Value = 0.02339786169299727
$ws.Range("J10").Value = 0.02339786169299728
$ws.Range("O10").Value = 0.01951796071871896
$ws.Range("P10").Value = 0.01951796071871896
$ws.Range("Q10").Value = 0.04157053910399999
$ws.Range("S10").Value = 0.0004566785454259397
$ws.Range("T10").Value = 0.0004566785454259399
